# Updated BOM for 3.81mm header.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 corresponds to P3, the 3.81mm vertical PCB header.
# Update its SuppliersPartNumber (column E) from C880557 to C8391,
# preserving the existing cell formatting (border + text style).
$target = $ws.Range("E3")
$target.Value = "C8391"

$formatSource = $ws.Range("D3")
$formatSource.Copy()
$target.PasteSpecial(-4122)  # xlPasteFormats

# Reflect the active selection left behind on save (cell E4).
$ws.Range("E4").Select()
